# Apply updated crypto price/volume data to the worksheet.
# A leading apostrophe forces Excel to store numeric-looking
# Price values (column D) as text, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.883.49'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.624.17'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '''211.61'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('D8').Value = '''23.02'
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('E10').Value = '  -1.32%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').Value = '1.854.97'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '1.625.62'
$ws.Range('E13').Value = '  -0.14%  '
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '''64.46'
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('D17').Value = '27.878.79'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '''227.63'
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0715'
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '''7.57'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').Value = '''9.92'
$ws.Range('E23').Value = '  -3.13%  '
$ws.Range('D24').Value = '''2.08'
$ws.Range('E24').Value = '  +2.16%  '
$ws.Range('D25').Value = '''154.12'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('D26').Value = '''6.90'
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('D32').Value = '''3.39'
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').Value = '1.420.50'
$ws.Range('E33').Value = '  +1.33%  '
$ws.Range('D34').Value = '''3.08'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('E35').Value = '  +1.76%  '
$ws.Range('D36').Value = '''0.982'
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('E39').Value = '  -0.64%  '
$ws.Range('E40').Value = '  -1.87%  '
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('D43').Value = '''65.24'
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('D45').Value = '''1.78'
$ws.Range('E45').Value = '  -3.18%  '
$ws.Range('D46').Value = '1.764.19'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('E47').Value = '  -4.00%  '
$ws.Range('D48').Value = '''89.18'
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('E51').Value = '  -0.37%  '
